# Commit: "Borrar filas innecesarias del Excel"
# A new method in the Excel DAL marks unneeded rows in the "Cheques" sheet
# with the text "BORRAR" (Spanish for "delete") so they can later be
# identified and removed. Mark the data rows of both the "Cheques" table
# (A8, A9) and the "Transferencias" table (A14, A15) accordingly, and leave
# the sheet's selection on C21, matching where the cursor ended up.

$wb = $excel.ActiveWorkbook
$wsCheques = $wb.Worksheets.Item("Cheques")

$wsCheques.Range("A8").Value = "BORRAR"
$wsCheques.Range("A9").Value = "BORRAR"
$wsCheques.Range("A14").Value = "BORRAR"
$wsCheques.Range("A15").Value = "BORRAR"

$wsCheques.Activate()
$wsCheques.Range("C21").Select()

# Restore the originally active sheet so the workbook-level active tab is
# unchanged by this edit.
$wb.Worksheets.Item("Resumen").Activate()
